$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("24.25 Budget Targets")

# Insert a new row above row 10 ("Arsenal v Dinamo Zagreb"); this shifts
# the existing rows 10-44 down to 11-45 and extends the sheet dimension.
$null = $ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new fixture.
$ws.Cells.Item(10, 1).Value = "UEFA Champions League"
$ws.Cells.Item(10, 2).Value = "Arsenal vs Real Madrid"
$ws.Cells.Item(10, 3).Value = 760000

# Grow Table1 so it covers the new row (A1:C44 -> A1:C45); this also keeps
# the AutoFilter range in sync.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:C45"))

# Keep the workbook-level _FilterDatabase defined name aligned with the
# table's new extent.
$fd = $wb.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "='24.25 Budget Targets'!`$A`$1:`$C`$45"

# Match the saved cursor position recorded in the workbook.
$null = $ws.Range("C11").Select()
